# Horarios actualizados Línea 141 - 200
# Applies the 07:16:53 scrape update to the three schedule sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Sheet "LP1912": refresh timestamp + swap two stop names back/forth ---
$ws1.Range("A2").Value = "Última actualización: 07:16:53"
$ws1.Range("C52").Value = "16_SANTA ANA"
$ws1.Range("C53").Value = "11_ETCHEVERRY"

# --- Sheet "LP1912-215": refresh timestamp only ---
$ws2.Range("A2").Value = "Última actualización: 07:16:53"

# --- Sheet "6203-6173": refresh timestamp, row count, and insert new scrape rows ---
$ws3.Range("A2").Value = "Última actualización: 07:16:53"
$ws3.Range("A3").Value = "Total filas: 16"

# New data, re-sorted ascending by Hora_Llegada (column B), rows 16-21
$data = @(
    @("07:16:53", "07:44", "215A_LA PLATA", 28,  "L6173"),
    @("06:23:52", "08:06", "215C_LA PLATA", 103, "L6203"),
    @("07:16:53", "08:10", "215C_LA PLATA", 54,  "L6203"),
    @("06:52:23", "08:11", "215C_LA PLATA", 79,  "L6203"),
    @("06:52:23", "08:40", "215A_LA PLATA", 108, "L6173"),
    @("07:16:53", "09:08", "215D_LA PLATA", 112, "L6203")
)

$row = 16
foreach ($entry in $data) {
    $ws3.Cells.Item($row, 1).Value = $entry[0]
    $ws3.Cells.Item($row, 2).Value = $entry[1]
    $ws3.Cells.Item($row, 3).Value = $entry[2]
    $ws3.Cells.Item($row, 4).Value = $entry[3]
    $ws3.Cells.Item($row, 5).Value = $entry[4]
    $row = $row + 1
}
